$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Fire Ninjutsu Scroll
$ws.Range("A2").Value = "Fire Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B2").Value = "MATK`n %`n1Staff only:Magic Pierce %`n5"
$ws.Range("C2").Value = "empty"
$ws.Range("D2").Value = "Sell1`nSpina`nProcess`n2 Wood"

# Row 3: Wind Ninjutsu Scroll
$ws.Range("A3").Value = "Wind Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B3").Value = "`nASPD250Katana only:`nCritical Rate`n5"
$ws.Range("C3").Value = "empty"
$ws.Range("D3").Value = "Sell1`nSpina`nProcess`n2 Wood"

# Row 4: Lightning Ninjutsu Scroll
$ws.Range("A4").Value = "Lightning Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B4").Value = "Stability %`n5Katana only:`nAccuracy %`n10"
$ws.Range("C4").Value = "empty"
$ws.Range("D4").Value = "Sell1`nSpina`nProcess`n2 Wood"

# Row 5: Metal Ninjutsu Scroll
$ws.Range("A5").Value = "Metal Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B5").Value = "`nCritical Rate`n5"
$ws.Range("C5").Value = "empty"
$ws.Range("D5").Value = "Sell1`nSpina`nProcess`n2 Wood"

# Row 6: Earth Ninjutsu Scroll
$ws.Range("A6").Value = "Earth Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B6").Value = "`nMaxHP`n %`n101-Handed Sword only:Fractional `nBarrier %`n10"
$ws.Range("C6").Value = "empty"
$ws.Range("D6").Value = "Sell1`nSpina`nProcess`n2 Wood"

# Row 7: Water Ninjutsu Scroll
$ws.Range("A7").Value = "Water Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B7").Value = "Ailment Resistance %`n5Magic Device only:`nAggro %`n-10"
$ws.Range("C7").Value = "empty"
$ws.Range("D7").Value = "Sell1`nSpina`nProcess`n2 Wood"

# Row 8: Dark Ninjutsu Scroll (unchanged content, only D column formatting changes)
$ws.Range("A8").Value = "Dark Ninjutsu Scroll[Ninjutsu Scroll]"
$ws.Range("B8").Value = "`nAggro %`n-10"
$ws.Range("C8").Value = "empty"
$ws.Range("D8").Value = "Sell1`nSpina`nProcess`n2 Wood"
